$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure all target cells keep their original text representation
# (avoids Excel auto-converting numeric-looking strings to floats,
# which would lose formatting like trailing zeros or thousands dots).

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '60.951.20'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +1.05%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.348.45'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -0.80%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '545.43'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +1.04%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '137.14'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -1.81%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.526'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -8.20%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.343.52'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -1.05%  '
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +0.27%  '
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +1.73%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '5.30'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -0.13%  '
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +0.71%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '24.69'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -2.05%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.773.74'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -0.88%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '60.921.43'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +1.06%  '
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -1.38%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.344.71'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -1.09%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '10.62'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +0.72%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '319.79'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +1.18%  '
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +1.19%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.55'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -1.45%  '
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +0.11%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.73'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -3.21%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '63.34'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +0.77%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.42'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +10.62%  '
$ws.Range('B27').NumberFormat = '@'
$ws.Range('B27').Value = 'WrappedeETH'
$ws.Range('C27').NumberFormat = '@'
$ws.Range('C27').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.464.04'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -1.00%  '
$ws.Range('B28').NumberFormat = '@'
$ws.Range('B28').Value = 'InternetComputer(DFINITY)'
$ws.Range('C28').NumberFormat = '@'
$ws.Range('C28').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.96'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +0.30%  '
$ws.Range('B29').NumberFormat = '@'
$ws.Range('B29').Value = 'Bittensor'
$ws.Range('C29').NumberFormat = '@'
$ws.Range('C29').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '500.59'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -2.40%  '
$ws.Range('B30').NumberFormat = '@'
$ws.Range('B30').Value = 'Fetch.AI'
$ws.Range('C30').NumberFormat = '@'
$ws.Range('C30').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.38'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -1.63%  '
$ws.Range('B31').NumberFormat = '@'
$ws.Range('B31').Value = 'PEPE'
$ws.Range('C31').NumberFormat = '@'
$ws.Range('C31').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0₃0866'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -5.86%  '
$ws.Range('B32').NumberFormat = '@'
$ws.Range('B32').Value = 'Kaspa'
$ws.Range('C32').NumberFormat = '@'
$ws.Range('C32').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.146'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +2.33%  '
$ws.Range('B33').NumberFormat = '@'
$ws.Range('B33').Value = 'PancakeSwap'
$ws.Range('C33').NumberFormat = '@'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.79'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -1.61%  '
$ws.Range('B34').NumberFormat = '@'
$ws.Range('B34').Value = 'ImmutableX'
$ws.Range('C34').NumberFormat = '@'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.50'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -3.43%  '
$ws.Range('B35').NumberFormat = '@'
$ws.Range('B35').Value = 'FirstDigitalUSD'
$ws.Range('C35').NumberFormat = '@'
$ws.Range('C35').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.00'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -0.03%  '
$ws.Range('B36').NumberFormat = '@'
$ws.Range('B36').Value = 'NEARProtocol'
$ws.Range('C36').NumberFormat = '@'
$ws.Range('C36').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.63'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +0.84%  '
$ws.Range('B37').NumberFormat = '@'
$ws.Range('B37').Value = 'PolygonEcosystemToken'
$ws.Range('C37').NumberFormat = '@'
$ws.Range('C37').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.377'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +1.43%  '
$ws.Range('B38').NumberFormat = '@'
$ws.Range('B38').Value = 'EthereumClassic'
$ws.Range('C38').NumberFormat = '@'
$ws.Range('C38').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '18.52'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +3.17%  '
$ws.Range('B39').NumberFormat = '@'
$ws.Range('B39').Value = 'Stacks'
$ws.Range('C39').NumberFormat = '@'
$ws.Range('C39').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.84'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +7.94%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.27'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -2.59%  '
$ws.Range('B41').NumberFormat = '@'
$ws.Range('B41').Value = 'Monero'
$ws.Range('C41').NumberFormat = '@'
$ws.Range('C41').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '140.93'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +3.58%  '
$ws.Range('B42').NumberFormat = '@'
$ws.Range('B42').Value = 'USDe'
$ws.Range('C42').NumberFormat = '@'
$ws.Range('C42').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.999'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -0.04%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '40.57'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +1.23%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '142.59'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +2.68%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.56'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +1.43%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.06'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -4.85%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0519'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +1.10%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '19.09'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -4.72%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.568'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -0.57%  '
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -1.72%  '
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -0.60%  '
